$wb = $excel.ActiveWorkbook

# Duplicate the last existing "Demand" sheet so the new sheet inherits
# identical formatting/styles, then rename and update its values.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "FTNC_Demand513"

# Update the data row values (header row stays identical to the template)
$newSheet.Range("B2").Value = 2201.685920603113
$newSheet.Range("C2").Value = 12909.84384652855
$newSheet.Range("D2").Value = 541.0858334763456
$newSheet.Range("E2").Value = 439.6554545412215
$newSheet.Range("F2").Value = 16092.27105514928
